# Update all the inventories, modify new_commodities in Add sector, run 2. Add sectors
#
# A new column D ("CF") is added to the Sheet1 path-inventory table, giving a
# third set of local folder paths (alongside the existing "LR" and "NG"
# columns) that point into Carol's local GitHub checkout.
#
# The cells are filled in the same order the original author appears to have
# used (top-to-bottom, skipping row 4 / "EXIOBASE Hybrid" which was left
# untouched, and coming back to fill D2 / MRSUT last).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("D1").Value = "CF"

# EXIOBASE IOT path (row 3) filled before going back for MRSUT (row 2)
$ws.Range("D3").Value = "C:\Users\carol\Desktop\UNI\MAGISTRALE\TESI\IOT"

# Row 4 (EXIOBASE Hybrid) is intentionally left untouched - no D4 value.

$ws.Range("D5").Value  = "C:\Users\carol\Desktop\UNI\MAGISTRALE\TESI\GitHub\GreenTechs\Database"
$ws.Range("D6").Value  = "C:\Users\carol\Desktop\UNI\MAGISTRALE\TESI\GitHub\GreenTechs\Add Sectors"
$ws.Range("D7").Value  = "C:\Users\carol\Desktop\UNI\MAGISTRALE\TESI\GitHub\GreenTechs\Shocks"
$ws.Range("D8").Value  = "C:\Users\carol\Desktop\UNI\MAGISTRALE\TESI\GitHub\GreenTechs\Results"
$ws.Range("D9").Value  = "C:\Users\carol\Desktop\UNI\MAGISTRALE\TESI\GitHub\GreenTechs\Plots"
$ws.Range("D10").Value = "C:\Users\carol\Desktop\UNI\MAGISTRALE\TESI\GitHub\GreenTechs\Shocks\ShockMaster.xlsx"

# EXIOBASE SUT / MRSUT path (row 2) filled in last
$ws.Range("D2").Value = "C:\Users\carol\Desktop\UNI\MAGISTRALE\TESI\MRSUT"

# Restore the selection the author left the sheet with
$ws.Range("G14").Select() | Out-Null
